# Rework "get data" header labels: strip the redundant hierarchical
# "Parent - Child - ..." prefixes from the row-1 header captions so each
# column just shows its own short label. Where stripping creates a
# collision with an already-used short label (the "Profits and dividends"
# Credit/Debit pair collides with the "Interests" Credit/Debit pair), the
# second occurrence is disambiguated with a ".2" suffix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @{
    "A1" = "datetime"
    "B1" = "Current Account"
    "C1" = "Goods (7)"
    "D1" = "Exports FOB"
    "E1" = "Imports FOB"
    "F1" = "Services"
    "G1" = "Exports"
    "H1" = "Imports"
    "I1" = "Income"
    "J1" = "Investment income"
    "K1" = "Interests"
    "L1" = "Credit"
    "M1" = "Debit"
    "N1" = "Profits and dividends"
    "O1" = "Credit.2"
    "P1" = "Debit.2"
    "Q1" = "Other Income"
    "R1" = "Current Transfers"
    "S1" = "Capital and Financial Account"
    "T1" = "Capital Account"
    "U1" = "Financial Account"
    "V1" = "Banking Sector"
    "W1" = "Central Bank"
    "X1" = "Other financial entities"
    "Y1" = "Nonfinancial Public Sector"
    "Z1" = "National Government (5) (6)"
    "AA1" = "Local Governments"
    "AB1" = "Companies and other"
    "AC1" = "Nonfinancial Private Sector"
    "AD1" = "Net Errors and Omissions"
    "AE1" = "International Reserves Variation"
    "AF1" = "BCRA International Reserves"
    "AG1" = "Exchange rate adjustment"
    "AH1" = "Imports CIF"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# The trailing row 24 was a stray, entirely empty row (only A24 had a
# style applied, no content) - remove it so the sheet's used range ends
# at row 23 again.
$ws.Rows.Item(24).Delete()

# Match the author's final active selection.
$ws.Range("P2").Select()
